# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns
# for the rows whose underlying data changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price (D) / Volume(1h) (E) cells are stored as plain text in the
# workbook (e.g. "299.28", "  -1.14%  "). The E cells always stay text
# automatically because of the leading/trailing spaces and "%" sign, but
# D cells whose text looks like a plain decimal number (e.g. "299.28")
# would otherwise get auto-converted into a real number by Excel. Force
# those specific D cells to Text format first so the refreshed values are
# written back as text, matching the original inline-string storage.
$textFormatCells = @("D6","D7","D8","D9","D10","D11","D12","D15","D16","D18","D19","D20","D21","D22","D23","D26","D27","D28","D29","D30","D32","D33","D34","D36","D38","D39","D40","D41","D42","D45","D46","D47","D48","D50","D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '23.450.32'
$ws.Range("E2").Value = '  -1.14%  '
$ws.Range("D3").Value = '1.646.96'
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("E5").Value = '  +0.21%  '
$ws.Range("D6").Value = '299.28'
$ws.Range("E6").Value = '  -1.60%  '
$ws.Range("D7").Value = '0.3790'
$ws.Range("E7").Value = '  -0.81%  '
$ws.Range("D8").Value = '0.3559'
$ws.Range("E8").Value = '  -1.38%  '
$ws.Range("D9").Value = '49.77'
$ws.Range("E9").Value = '  -3.28%  '
$ws.Range("D10").Value = '0.08097'
$ws.Range("E10").Value = '  -1.77%  '
$ws.Range("D11").Value = '1.220'
$ws.Range("E11").Value = '  -2.43%  '
$ws.Range("D12").Value = '1.002'
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("E13").Value = '  -2.67%  '
$ws.Range("E14").Value = '  -2.45%  '
$ws.Range("D15").Value = '7.364'
$ws.Range("E15").Value = '  -0.57%  '
$ws.Range("D16").Value = '0.00001195'
$ws.Range("D17").Value = '1.640.75'
$ws.Range("E17").Value = '  -0.60%  '
$ws.Range("D18").Value = '97.29'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("D19").Value = '0.06951'
$ws.Range("E19").Value = '  -0.30%  '
$ws.Range("D20").Value = '6.765'
$ws.Range("E20").Value = '  -0.35%  '
$ws.Range("D21").Value = '17.31'
$ws.Range("E21").Value = '  -2.34%  '
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("D23").Value = '12.40'
$ws.Range("D24").Value = '23.476.43'
$ws.Range("E24").Value = '  -1.05%  '
$ws.Range("E25").Value = '  -2.55%  '
$ws.Range("D26").Value = '2.930'
$ws.Range("E26").Value = '  -4.70%  '
$ws.Range("D27").Value = '20.87'
$ws.Range("E27").Value = '  -2.21%  '
$ws.Range("D28").Value = '152.83'
$ws.Range("E28").Value = '  +0.38%  '
$ws.Range("D29").Value = '5.211'
$ws.Range("E29").Value = '  -0.53%  '
$ws.Range("D30").Value = '132.65'
$ws.Range("E30").Value = '  -1.93%  '
$ws.Range("D31").Value = '1.830.68'
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").Value = '6.918'
$ws.Range("E32").Value = '  +0.55%  '
$ws.Range("D33").Value = '2.101'
$ws.Range("E33").Value = '  -0.25%  '
$ws.Range("D34").Value = '11.75'
$ws.Range("E34").Value = '  -1.68%  '
$ws.Range("E35").Value = '  -7.33%  '
$ws.Range("D36").Value = '0.02722'
$ws.Range("E36").Value = '  -3.86%  '
$ws.Range("E37").Value = '  -1.25%  '
$ws.Range("D38").Value = '0.2429'
$ws.Range("E38").Value = '  -3.65%  '
$ws.Range("D39").Value = '5.922'
$ws.Range("E39").Value = '  -2.80%  '
$ws.Range("D40").Value = '13.05'
$ws.Range("E40").Value = '  +1.71%  '
$ws.Range("D41").Value = '0.06775'
$ws.Range("E41").Value = '  -3.99%  '
$ws.Range("D42").Value = '0.6879'
$ws.Range("E43").Value = '  -2.78%  '
$ws.Range("E44").Value = '  -2.31%  '
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.17%  '
$ws.Range("D46").Value = '0.6394'
$ws.Range("E46").Value = '  -2.05%  '
$ws.Range("D47").Value = '2.251'
$ws.Range("E47").Value = '  -3.94%  '
$ws.Range("D48").Value = '3.915'
$ws.Range("E48").Value = '  -1.73%  '
$ws.Range("E49").Value = '  -3.42%  '
$ws.Range("D50").Value = '127.37'
$ws.Range("E50").Value = '  -0.61%  '
$ws.Range("D51").Value = '1.150'
$ws.Range("E51").Value = '  -3.41%  '
